$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 4933
$ws1.Range("F5").Value = 2833
$ws1.Range("F9").Value = 1744
$ws1.Range("F10").Value = 757
$ws1.Range("F11").Value = 496
$ws1.Range("F13").Value = 433
$ws1.Range("F14").Value = 1079
$ws1.Range("F15").Value = 305
$ws1.Range("F16").Value = 9
$ws1.Range("F19").Value = 1046
$ws1.Range("F20").Value = 47
$ws1.Range("F22").Value = 671
$ws1.Range("F23").Value = 755
$ws1.Range("F24").Value = 153
$ws1.Range("F25").Value = 12
$ws1.Range("F27").Value = 556
$ws1.Range("F28").Value = 53
$ws1.Range("F29").Value = 1663
$ws1.Range("F30").Value = 1656
$ws1.Range("F31").Value = 406
$ws1.Range("F33").Value = 1568
$ws1.Range("F34").Value = 220
$ws1.Range("F35").Value = 2399
$ws1.Range("F36").Value = 411
$ws1.Range("F37").Value = 29
$ws1.Range("F38").Value = 626
$ws1.Range("F39").Value = 117
$ws1.Range("F40").Value = 70
$ws1.Range("F42").Value = 819
$ws1.Range("F43").Value = 1511
$ws1.Range("F44").Value = 231
$ws1.Range("F46").Value = 503
$ws1.Range("F47").Value = 64
$ws1.Range("F48").Value = 83
$ws1.Range("F49").Value = 118

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 6
$ws2.Range("F4").Value = 105
$ws2.Range("F9").Value = 2
$ws2.Range("F12").Value = 47

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 4933
$ws4.Range("F4").Value = 2833
$ws4.Range("F5").Value = 1744
$ws4.Range("F7").Value = 6
$ws4.Range("F8").Value = 757
$ws4.Range("F9").Value = 496
$ws4.Range("F11").Value = 433
$ws4.Range("F12").Value = 1079
$ws4.Range("F13").Value = 305
$ws4.Range("F15").Value = 1046
$ws4.Range("F16").Value = 47
$ws4.Range("F17").Value = 671
$ws4.Range("F18").Value = 755
$ws4.Range("F19").Value = 153
$ws4.Range("F20").Value = 105
$ws4.Range("F21").Value = 105
$ws4.Range("F23").Value = 12
$ws4.Range("F26").Value = 556
$ws4.Range("F27").Value = 1663
$ws4.Range("F28").Value = 1656
$ws4.Range("F29").Value = 406
$ws4.Range("F33").Value = 2399
$ws4.Range("F34").Value = 411
$ws4.Range("F35").Value = 2
$ws4.Range("F38").Value = 29
$ws4.Range("F39").Value = 47
$ws4.Range("F40").Value = 117
$ws4.Range("F41").Value = 70
$ws4.Range("F43").Value = 819
$ws4.Range("F44").Value = 1511
$ws4.Range("F46").Value = 232
$ws4.Range("F47").Value = 503
$ws4.Range("F48").Value = 64
$ws4.Range("F49").Value = 83
